$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, shifting existing rows 31-68 down to 32-69
$ws.Rows.Item(31).Insert()

# Copy the style from the row above (row 30) for column D (date format) into new row 31
$ws.Cells.Item(30, 4).Copy()
$ws.Cells.Item(31, 4).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new row 31 with the new weekly data
$ws.Cells.Item(31, 1).Value = 2
$ws.Cells.Item(31, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(31, 3).Value = "Coquimbo"
$ws.Cells.Item(31, 4).Value = (Get-Date -Year 2022 -Month 11 -Day 30).Date
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(31, 6).Value = 100112032
$ws.Cells.Item(31, 7).Value = "Zapallo italiano"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 2400
$ws.Cells.Item(31, 11).Value = 5000
$ws.Cells.Item(31, 12).Value = 6000
$ws.Cells.Item(31, 13).Value = 5500
$ws.Cells.Item(31, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 92
$ws.Cells.Item(31, 17).Value = 60
$ws.Cells.Item(31, 18).Value = "Hortaliza"
